# Add the two new time-log entries (rows 50 and 51) that record putting the
# SkyScannerAPI idea on hold and starting the "manual input" window instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 50: 11/11/2023, 3 hours
$ws.Range("A50").Value2 = 45241
$ws.Range("A50").NumberFormat = "d-mmm"
$ws.Range("B50").Value2 = 3
$ws.Range("C50").Value2 = "so I tried to go back and use the other API, skyScanner API, so I went looking for the jar file since that was what was wrong.I couldn" + [char]0x2019 + "t find so I thought maybe I could add it as a depenedency, but no I crashed the while prpgram again, and had to re-install all the other HTTP jar files. InteleJ is not good"

# Row 51: 11/12/2023, 3 hours
$ws.Range("A51").Value2 = 45242
$ws.Range("A51").NumberFormat = "d-mmm"
$ws.Range("B51").Value2 = 3
$ws.Range("C51").Value2 = "I put the dream of having a flight api aside, and Im just going to have someone input their data into the a little window that" + [char]0x2019 + "s going to popup, I added the code to my controller class"

# Match the author's final selection / scroll position from the saved file.
$ws.Range("A36").Select()
$excel.ActiveWindow.ScrollRow = 36
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E54").Select()
